$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.608.16"
$ws.Range("E2").Value = "  +1.06%  "
$ws.Range("D3").Value = "3.392.46"
$ws.Range("E3").Value = "  +0.20%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'576.75"
$ws.Range("E5").Value = "  +0.89%  "
$ws.Range("D6").Value = "'142.95"
$ws.Range("E6").Value = "  +0.53%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -0.28%  "
$ws.Range("E9").Value = "  -0.08%  "
$ws.Range("D10").Value = "'0.122"
$ws.Range("E10").Value = "  -0.53%  "
$ws.Range("E11").Value = "  -1.01%  "
$ws.Range("D12").Value = "3.972.62"
$ws.Range("E12").Value = "  +0.21%  "
$ws.Range("E13").Value = "  -0.24%  "
$ws.Range("D14").Value = "'27.97"
$ws.Range("E14").Value = "  +0.72%  "
$ws.Range("D15").Value = "3.395.14"
$ws.Range("E15").Value = "  +0.77%  "
$ws.Range("E16").Value = "  -0.96%  "
$ws.Range("D17").Value = "61.623.65"
$ws.Range("E17").Value = "  +0.91%  "
$ws.Range("E18").Value = "  +0.61%  "
$ws.Range("D19").Value = "'13.71"
$ws.Range("E19").Value = "  +0.31%  "
$ws.Range("D20").Value = "'9.13"
$ws.Range("E20").Value = "  +1.53%  "
$ws.Range("D21").Value = "'387.55"
$ws.Range("E21").Value = "  +0.84%  "
$ws.Range("D22").Value = "'74.49"
$ws.Range("E22").Value = "  -0.62%  "
$ws.Range("E23").Value = "  -0.74%  "
$ws.Range("E24").Value = "  -0.14%  "
$ws.Range("D25").Value = "'0.0000113"
$ws.Range("E25").Value = "  -2.68%  "
$ws.Range("E26").Value = "  -0.50%  "
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("D28").Value = "'7.37"
$ws.Range("E28").Value = "  +1.10%  "
$ws.Range("E29").Value = "  -0.32%  "
$ws.Range("E30").Value = "  -0.54%  "
$ws.Range("B31").Value = "USDe"
$ws.Range("C31").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").Value = "'1.39"
$ws.Range("E32").Value = "  -0.22%  "
$ws.Range("D33").Value = "'23.35"
$ws.Range("E33").Value = "  +0.21%  "
$ws.Range("D34").Value = "'6.92"
$ws.Range("E34").Value = "  -0.54%  "
$ws.Range("D35").Value = "'168.85"
$ws.Range("E35").Value = "  +1.17%  "
$ws.Range("D36").Value = "'5.12"
$ws.Range("E36").Value = "  +2.29%  "
$ws.Range("D37").Value = "3.424.19"
$ws.Range("E37").Value = "  +0.18%  "
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("D39").Value = "'27.35"
$ws.Range("E39").Value = "  +2.69%  "
$ws.Range("D40").Value = "'0.0757"
$ws.Range("E40").Value = "  -1.34%  "
$ws.Range("D41").Value = "'0.782"
$ws.Range("E41").Value = "  +0.40%  "
$ws.Range("E42").Value = "  +0.85%  "
$ws.Range("D43").Value = "'1.67"
$ws.Range("E43").Value = "  -0.20%  "
$ws.Range("E44").Value = "  +2.52%  "
$ws.Range("D45").Value = "2.474.82"
$ws.Range("E45").Value = "  +0.56%  "
$ws.Range("D46").Value = "'22.73"
$ws.Range("E46").Value = "  -1.40%  "
$ws.Range("D47").Value = "'6.63"
$ws.Range("E47").Value = "  -1.61%  "
$ws.Range("D48").Value = "'0.999"
$ws.Range("E48").Value = "  -0.02%  "
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("E50").Value = "  -6.08%  "
$ws.Range("E51").Value = "  -1.47%  "
